$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D5,D6,D8,D13,D14,D15,D19,D20,D21,D22,D23,D24,D25,D26,D27,D28,D29,D31,D32,D33,D34,D35,D36,D37,D38,D39,D40,D41,D42,D43,D44,D45,D46,D47,D48,D49,D50,D51').NumberFormat = "@"

$ws.Range('D2').Value = '61.463.99'
$ws.Range('E2').Value = '  +1.38%  '

$ws.Range('D3').Value = '2.380.88'
$ws.Range('E3').Value = '  +1.52%  '

$ws.Range('E4').Value = '  +0.17%  '

$ws.Range('D5').Value = '552.42'
$ws.Range('E5').Value = '  +1.79%  '

$ws.Range('D6').Value = '139.96'
$ws.Range('E6').Value = '  +1.91%  '

$ws.Range('E7').Value = '  +0.11%  '

$ws.Range('D8').Value = '0.524'
$ws.Range('E8').Value = '  +1.45%  '

$ws.Range('D9').Value = '2.382.07'
$ws.Range('E9').Value = '  +1.57%  '

$ws.Range('E10').Value = '  +5.03%  '

$ws.Range('E11').Value = '  +2.06%  '

$ws.Range('E12').Value = '  +3.21%  '

$ws.Range('D13').Value = '0.352'
$ws.Range('E13').Value = '  +4.43%  '

$ws.Range('D14').Value = '25.54'
$ws.Range('E14').Value = '  +4.32%  '

$ws.Range('D15').Value = '0.0000168'
$ws.Range('E15').Value = '  +5.62%  '

$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '2.813.05'
$ws.Range('E16').Value = '  +1.82%  '

$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '61.206.62'
$ws.Range('E17').Value = '  +1.49%  '

$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.381.85'
$ws.Range('E18').Value = '  +1.62%  '

$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '11.03'
$ws.Range('E19').Value = '  +5.08%  '

$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = '4.16'
$ws.Range('E20').Value = '  +2.69%  '

$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '321.58'
$ws.Range('E21').Value = '  +3.17%  '

$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '6.76'
$ws.Range('E22').Value = '  +4.05%  '

$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.20%  '

$ws.Range('B24').Value = 'SuiNetwork'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D24').Value = '1.78'
$ws.Range('E24').Value = '  -4.57%  '

$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '64.44'
$ws.Range('E25').Value = '  +2.98%  '

$ws.Range('B26').Value = 'Aptos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D26').Value = '8.93'
$ws.Range('E26').Value = '  +10.84%  '

$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.09%  '

$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '8.22'
$ws.Range('E28').Value = '  +4.01%  '

$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').Value = '518.71'
$ws.Range('E29').Value = '  +3.27%  '

$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0903'
$ws.Range('E30').Value = '  +2.28%  '

$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').Value = '0.151'
$ws.Range('E31').Value = '  +5.13%  '

$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').Value = '1.40'
$ws.Range('E32').Value = '  +1.41%  '

$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '1.85'
$ws.Range('E33').Value = '  +3.60%  '

$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '1.55'
$ws.Range('E34').Value = '  +1.60%  '

$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.08%  '

$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D36').Value = '5.56'
$ws.Range('E36').Value = '  +6.14%  '

$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').Value = '4.72'
$ws.Range('E37').Value = '  +4.77%  '

$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '1.89'
$ws.Range('E38').Value = '  +6.12%  '

$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D39').Value = '0.379'
$ws.Range('E39').Value = '  +2.50%  '

$ws.Range('B40').Value = 'EthereumClassic'
$ws.Range('C40').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D40').Value = '18.54'
$ws.Range('E40').Value = '  +0.75%  '

$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').Value = '146.61'
$ws.Range('E41').Value = '  +6.27%  '

$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.06%  '

$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').Value = '41.32'
$ws.Range('E43').Value = '  +2.99%  '

$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '148.37'
$ws.Range('E44').Value = '  +7.84%  '

$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = '2.16'
$ws.Range('E45').Value = '  +4.91%  '

$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').Value = '3.61'
$ws.Range('E46').Value = '  +2.34%  '

$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').Value = '0.0527'
$ws.Range('E47').Value = '  +3.95%  '

$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '19.76'
$ws.Range('E48').Value = '  +2.38%  '

$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '0.582'
$ws.Range('E49').Value = '  +3.15%  '

$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').Value = '0.0906'
$ws.Range('E50').Value = '  +1.69%  '

$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').Value = '0.0225'
$ws.Range('E51').Value = '  +2.10%  '

Write-Host "Updated cryptos list"